# Add a "TotalUse" column to the Results sheet, between "MainUse" and
# "Result", containing the sum of the four StackUse columns plus MainUse,
# and hook it into Table1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")
$lo = $ws.ListObjects.Item("Table1")

$lastRow = $lo.Range.Rows.Count + $lo.Range.Row - 1   # 41

# Insert a blank worksheet column at G (this shifts the existing "Result"
# table column, and its data/formatting, one place to the right, to H).
$ws.Range("G1").EntireColumn.Insert()

# Grow the table definition to cover the freshly inserted column.
$lo.Resize($ws.Range("A1:H" + $lastRow))

# Name the header cells -- this also renames the matching table columns.
$ws.Range("G1").Value2 = "TotalUse"
$ws.Range("H1").Value2 = "Result"

# Fill in the calculated values for the new column, row by row, so that
# each row gets its own formula (matching how the table's calculated
# column behaves).
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Formula = "=SUM(Table1[[#This Row],[StackUse1]:[MainUse]])"
}

# Match the column width used for the rest of the numeric columns.
$ws.Range("G1").EntireColumn.ColumnWidth = $ws.Range("F1").EntireColumn.ColumnWidth

# Reflect the selection left behind by the edit.
$ws.Range("H6").Select()
